# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col D) and
# "Correspond Handback DateTime" (col G) timestamps for rows 2 and 3
# on the "zh-cn" and "de-de" worksheets.

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D2").Value = "2016-02-23 09:29:51"
$wsZh.Range("D3").Value = "2016-02-23 09:29:51"
$wsZh.Range("G2").Value = "2016-02-23 09:31:14"
$wsZh.Range("G3").Value = "2016-02-23 09:31:14"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D2").Value = "2016-02-23 09:30:02"
$wsDe.Range("D3").Value = "2016-02-23 09:30:02"
$wsDe.Range("G2").Value = "2016-02-23 09:31:46"
$wsDe.Range("G3").Value = "2016-02-23 09:31:46"
